# Auto-generated Excel COM-interop script
# Commit: "case with 380 kV done"
# Updates numeric result values in Sheet1 (res_line/pl_mw data), rows 2-25,
# columns B, C, E, F, G, H, I, J, K, L, O, to reflect the re-run with a 380 kV case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.5404584315820387
$ws.Range("C2").Value = 0.2338570918008038
$ws.Range("E2").Value = 0.2658227537422171
$ws.Range("F2").Value = 1.796121191362907
$ws.Range("G2").Value = 0.5927931663030961
$ws.Range("H2").Value = 0.736074104013241
$ws.Range("I2").Value = 0.7865246678824036
$ws.Range("J2").Value = 0.04566136993263825
$ws.Range("K2").Value = 0.2955353816041963
$ws.Range("L2").Value = 0.5029601225827207
$ws.Range("O2").Value = 2.638485371817936
# Row 3
$ws.Range("B3").Value = 0.4934470909824711
$ws.Range("C3").Value = 0.234980377440607
$ws.Range("E3").Value = 0.264449994233324
$ws.Range("F3").Value = 1.798864819960549
$ws.Range("G3").Value = 0.5991942512098518
$ws.Range("H3").Value = 0.7425893750361894
$ws.Range("I3").Value = 0.7951521567995421
$ws.Range("J3").Value = 0.04367312395655532
$ws.Range("K3").Value = 0.2582842804359302
$ws.Range("L3").Value = 0.4911736614751163
$ws.Range("O3").Value = 2.66547099865609
# Row 4
$ws.Range("B4").Value = 0.4646041524358679
$ws.Range("C4").Value = 0.2357164922915125
$ws.Range("E4").Value = 0.2637085531744425
$ws.Range("F4").Value = 1.801453345606362
$ws.Range("G4").Value = 0.6034952047385218
$ws.Range("H4").Value = 0.7468777035591785
$ws.Range("I4").Value = 0.8008285796039072
$ws.Range("J4").Value = 0.04244396534503636
$ws.Range("K4").Value = 0.2353318005448699
$ws.Range("L4").Value = 0.4841208151304386
$ws.Range("O4").Value = 2.683421191825957
# Row 5
$ws.Range("B5").Value = 0.4528569251340286
$ws.Range("C5").Value = 0.2360281654443988
$ws.Range("E5").Value = 0.2634319984646503
$ws.Range("F5").Value = 1.802735752769713
$ws.Range("G5").Value = 0.6053410472298566
$ws.Range("H5").Value = 0.7486976890046719
$ws.Range("I5").Value = 0.8032371182609594
$ws.Range("J5").Value = 0.04194099526712591
$ws.Range("K5").Value = 0.2259589925445198
$ws.Range("L5").Value = 0.481293294268113
$ws.Range("O5").Value = 2.691083315074465
# Row 6
$ws.Range("B6").Value = 0.4509067268996318
$ws.Range("C6").Value = 0.2360806260781452
$ws.Range("E6").Value = 0.263387624629047
$ws.Range("F6").Value = 1.802962447383088
$ws.Range("G6").Value = 0.6056531740641802
$ws.Range("H6").Value = 0.7490042740437985
$ws.Range("I6").Value = 0.8036428137820959
$ws.Range("J6").Value = 0.04185735273294355
$ws.Range("K6").Value = 0.2244014877364862
$ws.Range("L6").Value = 0.4808266076321956
$ws.Range("O6").Value = 2.692376580255484
# Row 7
$ws.Range("B7").Value = 0.4644456976656102
$ws.Range("C7").Value = 0.2357206482109113
$ws.Range("E7").Value = 0.2637047197463431
$ws.Range("F7").Value = 1.801469718860019
$ws.Range("G7").Value = 0.6035197212269665
$ws.Range("H7").Value = 0.7469019550709959
$ws.Range("I7").Value = 0.8008606759021646
$ws.Range("J7").Value = 0.04243719049367201
$ws.Range("K7").Value = 0.2352054735881524
$ws.Range("L7").Value = 0.4840824932567784
$ws.Range("O7").Value = 2.683523119757325
# Row 8
$ws.Range("B8").Value = 0.5242448462121558
$ws.Range("C8").Value = 0.2342347853734239
$ws.Range("E8").Value = 0.2653284202461101
$ws.Range("F8").Value = 1.796879738327
$ws.Range("G8").Value = 0.5949233076132785
$ws.Range("H8").Value = 0.7382608441163612
$ws.Range("I8").Value = 0.7894207544706902
$ws.Range("J8").Value = 0.04497757453510332
$ws.Range("K8").Value = 0.2827082273123551
$ws.Range("L8").Value = 0.4988581071498714
$ws.Range("O8").Value = 2.647503429053202
# Row 9
$ws.Range("B9").Value = 0.641652501972743
$ws.Range("C9").Value = 0.2316879629434752
$ws.Range("E9").Value = 0.2693142003464324
$ws.Range("F9").Value = 1.795041268866356
$ws.Range("G9").Value = 0.5810085238595875
$ws.Range("H9").Value = 0.7235976707901273
$ws.Range("I9").Value = 0.7699941434043467
$ws.Range("J9").Value = 0.04989193515774559
$ws.Range("K9").Value = 0.3751999940202779
$ws.Range("L9").Value = 0.5292831937632201
$ws.Range("O9").Value = 2.587824384661872
# Row 10
$ws.Range("B10").Value = 0.7279603686359053
$ws.Range("C10").Value = 0.2300387145516538
$ws.Range("E10").Value = 0.2727276217359318
$ws.Range("F10").Value = 1.798046397050697
$ws.Range("G10").Value = 0.5725817943414384
$ws.Range("H10").Value = 0.714212243552808
$ws.Range("I10").Value = 0.7575529944186243
$ws.Range("J10").Value = 0.05346053797094896
$ws.Range("K10").Value = 0.4427235789582085
$ws.Range("L10").Value = 0.5525091298193701
$ws.Range("O10").Value = 2.55065503118405
# Row 11
$ws.Range("B11").Value = 0.7672270711624662
$ws.Range("C11").Value = 0.2293362340898355
$ws.Range("E11").Value = 0.2743850592222898
$ws.Range("F11").Value = 1.80035725314761
$ws.Range("G11").Value = 0.5691388836356168
$ws.Range("H11").Value = 0.7102431209998343
$ws.Range("I11").Value = 0.752290567675459
$ws.Range("J11").Value = 0.05507469140405874
$ws.Range("K11").Value = 0.4733430656582698
$ws.Range("L11").Value = 0.5632624270363067
$ws.Range("O11").Value = 2.535195181888227
# Row 12
$ws.Range("B12").Value = 0.7820962251465176
$ws.Range("C12").Value = 0.2290770633159092
$ws.Range("E12").Value = 0.2750276627524784
$ws.Range("F12").Value = 1.801367809637938
$ws.Range("G12").Value = 0.5678913357240774
$ws.Range("H12").Value = 0.7087832551662956
$ws.Range("I12").Value = 0.7503549119724653
$ws.Range("J12").Value = 0.05568458018851885
$ws.Range("K12").Value = 0.4849232811754405
$ws.Range("L12").Value = 0.5673611623258097
$ws.Range("O12").Value = 2.52954924948736
# Row 13
$ws.Range("B13").Value = 0.7788939155636854
$ws.Range("C13").Value = 0.2291325764317236
$ws.Range("E13").Value = 0.2748886021785353
$ws.Range("F13").Value = 1.801144146195355
$ws.Range("G13").Value = 0.5681575166691033
$ws.Range("H13").Value = 0.7090957446459427
$ws.Range("I13").Value = 0.7507692499267655
$ws.Range("J13").Value = 0.05555329051784952
$ws.Range("K13").Value = 0.4824299407169121
$ws.Range("L13").Value = 0.5664772438001933
$ws.Range("O13").Value = 2.530755933493623
# Row 14
$ws.Range("B14").Value = 0.7684503769496018
$ws.Range("C14").Value = 0.2293147749534512
$ws.Range("E14").Value = 0.2744376270881759
$ws.Range("F14").Value = 1.800437678249295
$ws.Range("G14").Value = 0.5690351202841555
$ws.Range("H14").Value = 0.7101221522796379
$ws.Range("I14").Value = 0.7521301754870997
$ws.Range("J14").Value = 0.05512489467150772
$ws.Range("K14").Value = 0.4742960760956976
$ws.Range("L14").Value = 0.5635990991467139
$ws.Range("O14").Value = 2.534726511110264
# Row 15
$ws.Range("B15").Value = 0.7620533404433445
$ws.Range("C15").Value = 0.2294272672425457
$ws.Range("E15").Value = 0.2741633386737874
$ws.Range("F15").Value = 1.800022583490929
$ws.Range("G15").Value = 0.5695799995418938
$ws.Range("H15").Value = 0.7107564755483224
$ws.Range("I15").Value = 0.752971219432407
$ws.Range("J15").Value = 0.05486231235155969
$ws.Range("K15").Value = 0.4693119108024746
$ws.Range("L15").Value = 0.5618396201183771
$ws.Range("O15").Value = 2.53718574312623
# Row 16
$ws.Range("B16").Value = 0.7253942098307675
$ws.Range("C16").Value = 0.2300855823399708
$ws.Range("E16").Value = 0.2726214032739165
$ws.Range("F16").Value = 1.797914349124071
$ws.Range("G16").Value = 0.5728146587771974
$ws.Range("H16").Value = 0.7144776758801967
$ws.Range("I16").Value = 0.7579048978978307
$ws.Range("J16").Value = 0.05335486111802368
$ws.Range("K16").Value = 0.4407205030412626
$ws.Range("L16").Value = 0.5518101282824119
$ws.Range("O16").Value = 2.551694524836165
# Row 17
$ws.Range("B17").Value = 0.7029055520346219
$ws.Range("C17").Value = 0.2305016539628575
$ws.Range("E17").Value = 0.27170222266313
$ws.Range("F17").Value = 1.796862582296754
$ws.Range("G17").Value = 0.5748990650962824
$ws.Range("H17").Value = 0.7168374142684542
$ws.Range("I17").Value = 0.7610332650064997
$ws.Range("J17").Value = 0.05242770535426899
$ws.Range("K17").Value = 0.4231551561044569
$ws.Range("L17").Value = 0.5457052242169169
$ws.Range("O17").Value = 2.560966275531086
# Row 18
$ws.Range("B18").Value = 0.6899711835039
$ws.Range("C18").Value = 0.23074546531927
$ws.Range("E18").Value = 0.2711833897787201
$ws.Range("F18").Value = 1.796346483575192
$ws.Range("G18").Value = 0.5761347079369017
$ws.Range("H18").Value = 0.7182229463798109
$ws.Range("I18").Value = 0.762869996916983
$ws.Range("J18").Value = 0.05189356367232278
$ws.Range("K18").Value = 0.4130429149853114
$ws.Range("L18").Value = 0.5422115300446961
$ws.Range("O18").Value = 2.566435494063228
# Row 19
$ws.Range("B19").Value = 0.6855919467358831
$ws.Range("C19").Value = 0.2308287889836471
$ws.Range("E19").Value = 0.271009417045228
$ws.Range("F19").Value = 1.796187008516881
$ws.Range("G19").Value = 0.5765593849042077
$ws.Range("H19").Value = 0.7186969207634561
$ws.Range("I19").Value = 0.7634983022721542
$ws.Range("J19").Value = 0.0517125646170129
$ws.Range("K19").Value = 0.4096175399323556
$ws.Range("L19").Value = 0.5410316716993009
$ws.Range("O19").Value = 2.568310694388288
# Row 20
$ws.Range("B20").Value = 0.7052994621531354
$ws.Range("C20").Value = 0.2304568970958272
$ws.Range("E20").Value = 0.2717990515910031
$ws.Range("F20").Value = 1.796965350739669
$ws.Range("G20").Value = 0.5746733728152762
$ws.Range("H20").Value = 0.7165832902002407
$ws.Range("I20").Value = 0.760696376468033
$ws.Range("J20").Value = 0.05252649253451835
$ws.Range("K20").Value = 0.4250259662146334
$ws.Range("L20").Value = 0.5463532734734002
$ws.Range("O20").Value = 2.559965169410134
# Row 21
$ws.Range("B21").Value = 0.7715179132321168
$ws.Range("C21").Value = 0.2292610733113136
$ws.Range("E21").Value = 0.2745696838275862
$ws.Range("F21").Value = 1.800641509840233
$ws.Range("G21").Value = 0.5687758208844969
$ws.Range("H21").Value = 0.7098195005554473
$ws.Range("I21").Value = 0.7517288886301721
$ws.Range("J21").Value = 0.05525076195937118
$ws.Range("K21").Value = 0.4766855945414363
$ws.Range("L21").Value = 0.564443757867366
$ws.Range("O21").Value = 2.533554600457393
# Row 22
$ws.Range("B22").Value = 0.8147933286037983
$ws.Range("C22").Value = 0.228519409610616
$ws.Range("E22").Value = 0.276467658860156
$ws.Range("F22").Value = 1.803833695753042
$ws.Range("G22").Value = 0.5652490789504725
$ws.Range("H22").Value = 0.7056504953312128
$ws.Range("I22").Value = 0.7462010071032488
$ws.Range("J22").Value = 0.05702331120844661
$ws.Range("K22").Value = 0.5103620171088039
$ws.Range("L22").Value = 0.576422396373971
$ws.Range("O22").Value = 2.517508374364766
# Row 23
$ws.Range("B23").Value = 0.7916969491094505
$ws.Range("C23").Value = 0.2289116093416581
$ws.Range("E23").Value = 0.2754467213474072
$ws.Range("F23").Value = 1.802057794739298
$ws.Range("G23").Value = 0.5671013674927252
$ws.Range("H23").Value = 0.7078525693966355
$ws.Range("I23").Value = 0.7491208798883413
$ws.Range("J23").Value = 0.05607800277320507
$ws.Range("K23").Value = 0.4923963959838318
$ws.Range("L23").Value = 0.5700150448643342
$ws.Range("O23").Value = 2.525961391244948
# Row 24
$ws.Range("B24").Value = 0.7042171919457019
$ws.Range("C24").Value = 0.2304771173248028
$ws.Range("E24").Value = 0.2717552452733116
$ws.Range("F24").Value = 1.796918613203758
$ws.Range("G24").Value = 0.5747752921453895
$ws.Range("H24").Value = 0.7166980897076414
$ws.Range("I24").Value = 0.7608485648117522
$ws.Range("J24").Value = 0.05248183429780795
$ws.Range("K24").Value = 0.4241802155493133
$ws.Range("L24").Value = 0.5460602402444721
$ws.Range("O24").Value = 2.5604173372476
# Row 25
$ws.Range("B25").Value = 0.6098793810860457
$ws.Range("C25").Value = 0.2323378486123531
$ws.Range("E25").Value = 0.2681504923423397
$ws.Range("F25").Value = 1.794772920484803
$ws.Range("G25").Value = 0.5844575410406492
$ws.Range("H25").Value = 0.7273205125535398
$ws.Range("I25").Value = 0.7749277967386128
$ws.Range("J25").Value = 0.04856977431451526
$ws.Range("K25").Value = 0.3502521017344407
$ws.Range("L25").Value = 0.520898333285345
$ws.Range("O25").Value = 2.602796485132856

Write-Output "Updated B2:O25 (excl. D, M, N) values for 380 kV case"
